# Add a "decision tree regressor" row to the scores sheet, matching the
# mse / rmse columns already present for the other two methods, then
# reformat the numeric score columns as percentages and widen column A
# to fit the new (longer) label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 4): method name + its two metrics.
$ws.Range("A4").Value = "decision tree regressor"
$ws.Range("B4").Value = 0.22163588390501299
$ws.Range("C4").Value = 0.47078220432065299

# Reformat the metric columns (now B2:C4) as percentages with 2 decimals.
$ws.Range("B2:C4").NumberFormat = "0.00%"

# Widen column A so the longer "decision tree regressor" label fits
# (matches a displayed character width of 21).
$ws.Columns.Item(1).ColumnWidth = 20.166666666666668

# Leave the selection where the author ended up after editing.
$ws.Range("C11").Select() | Out-Null
